$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28
$ws.Range("A28").Value = 131089521
$ws.Range("B28").Value = 57881
$ws.Range("E28").Value = 100049
$ws.Range("F28").Value = "Spillkråka"
$ws.Range("G28").Value = "Dryocopus martius"
$ws.Range("L28").ClearContents()
$ws.Range("M28").Value = "färska spår"
$ws.Range("Q28").Value = 584995
$ws.Range("R28").Value = 7060537
$ws.Range("Z28").NumberFormat = "@"
$ws.Range("Z28").Value = "13:14"
$ws.Range("AB28").NumberFormat = "@"
$ws.Range("AB28").Value = "13:14"

# Row 29
$ws.Range("A29").Value = 131090275
$ws.Range("B29").Value = 57884
$ws.Range("E29").Value = 100109
$ws.Range("F29").Value = "Tretåig hackspett"
$ws.Range("G29").Value = "Picoides tridactylus"
$ws.Range("L29").Value = "hane"
$ws.Range("M29").Value = "födosökande"
$ws.Range("Q29").Value = 584987
$ws.Range("R29").Value = 7060190
$ws.Range("Z29").ClearContents()
$ws.Range("AB29").ClearContents()

# Row 32
$ws.Range("A32").Value = 131144496
$ws.Range("B32").Value = 57884
$ws.Range("E32").Value = 100109
$ws.Range("F32").Value = "Tretåig hackspett"
$ws.Range("G32").Value = "Picoides tridactylus"
$ws.Range("H32").Value = "(Linnaeus, 1758)"
$ws.Range("M32").Value = "färska spår"
$ws.Range("Q32").Value = 584875
$ws.Range("R32").Value = 7060422
$ws.Range("AC32").Value = "Färska ringhack, tall"

# Row 33
$ws.Range("A33").Value = 131144498
$ws.Range("B33").Value = 58043
$ws.Range("E33").Value = 103021
$ws.Range("F33").Value = "Talltita"
$ws.Range("G33").Value = "Poecile montanus"
$ws.Range("H33").Value = "(Conrad von Baldenstein, 1827)"
$ws.Range("M33").Value = "lockläte, övriga läten"
$ws.Range("Q33").Value = 584857
$ws.Range("R33").Value = 7060494
$ws.Range("AC33").ClearContents()
